$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 36 for the new UPGD "IPS CENTRO DE MEDICINA INTEGRATIVA SAS" (shifts rows 36-52 down to 37-53)
$ws.Rows.Item(36).Insert()

# Fill in the data for the newly inserted row 36
$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = "6600102402"
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "03"
$ws.Range("C36").Value = "IPS CENTRO DE MEDICINA INTEGRATIVA SAS"
$ws.Range("Q36").Value = 0
$ws.Range("S36").Value = 0

# Rename UPGD in row 29 (LIGA CONTRA EL CANCER SECCIONAL RISARALDA -> FUNDACION LA LIGA AMA SALVAR VIDAS)
$ws.Range("C29").Value = "FUNDACIÓN LA LIGA AMA SALVAR VIDAS"

# Correct week-12/13 values on row 26 (Clinica Los Rosales): O26 2->10, P26 3->2
$ws.Range("O26").Value = 10
$ws.Range("P26").Value = 2

# Add new week columns (14, 15, 16 -> Q, R, S) data for semana 16 de 2025
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("Q6").Value = 15
$ws.Range("R6").Value = 20
$ws.Range("S6").Value = 16
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = 2
$ws.Range("S7").Value = 1
$ws.Range("Q8").Value = 10
$ws.Range("R8").Value = 18
$ws.Range("S8").Value = 28
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("R14").Value = 0
$ws.Range("R15").Value = 0
$ws.Range("S15").Value = 0
$ws.Range("R16").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("Q21").Value = 0
$ws.Range("R21").Value = 0
$ws.Range("Q23").Value = 3
$ws.Range("R23").Value = 1
$ws.Range("S23").Value = 1
$ws.Range("Q24").Value = 0
$ws.Range("R24").Value = 0
$ws.Range("S24").Value = 0
$ws.Range("R25").Value = 0
$ws.Range("Q26").Value = 2
$ws.Range("S26").Value = 0
$ws.Range("Q27").Value = 3
$ws.Range("R27").Value = 2
$ws.Range("S27").Value = 14
$ws.Range("Q28").Value = 17
$ws.Range("R28").Value = 14
$ws.Range("S28").Value = 20
$ws.Range("Q29").Value = 0
$ws.Range("R29").Value = 0
$ws.Range("S29").Value = 0
$ws.Range("Q30").Value = 0
$ws.Range("R30").Value = 0
$ws.Range("S30").Value = 0
$ws.Range("R31").Value = 0
$ws.Range("Q32").Value = 13
$ws.Range("R32").Value = 26
$ws.Range("S32").Value = 16
$ws.Range("Q33").Value = 4
$ws.Range("R33").Value = 0
$ws.Range("S33").Value = 4
$ws.Range("Q34").Value = 0
$ws.Range("R34").Value = 0
$ws.Range("Q35").Value = 0
$ws.Range("R35").Value = 0
$ws.Range("S35").Value = 0
$ws.Range("Q37").Value = 0
$ws.Range("R37").Value = 0
$ws.Range("S37").Value = 0
$ws.Range("Q38").Value = 0
$ws.Range("R38").Value = 0
$ws.Range("S38").Value = 0
$ws.Range("Q39").Value = 0
$ws.Range("R39").Value = 0
$ws.Range("S39").Value = 0
$ws.Range("Q40").Value = 0
$ws.Range("R40").Value = 0
$ws.Range("S40").Value = 0
$ws.Range("Q41").Value = 0
$ws.Range("R41").Value = 0
$ws.Range("S41").Value = 0
$ws.Range("Q42").Value = 0
$ws.Range("R42").Value = 0
$ws.Range("S42").Value = 0
$ws.Range("Q43").Value = 0
$ws.Range("R43").Value = 0
$ws.Range("S43").Value = 0
$ws.Range("Q44").Value = 0
$ws.Range("R44").Value = 0
$ws.Range("S44").Value = 0
$ws.Range("Q45").Value = 0
$ws.Range("R45").Value = 0
$ws.Range("S45").Value = 0
$ws.Range("Q46").Value = 0
$ws.Range("R46").Value = 0
$ws.Range("S46").Value = 0
$ws.Range("Q47").Value = 0
$ws.Range("S47").Value = 0
$ws.Range("Q48").Value = 0
$ws.Range("R48").Value = 0
$ws.Range("Q49").Value = 0
$ws.Range("R49").Value = 0
$ws.Range("S49").Value = 1
$ws.Range("Q50").Value = 0
$ws.Range("R50").Value = 0
$ws.Range("S50").Value = 0
$ws.Range("Q51").Value = 0
$ws.Range("R51").Value = 0
$ws.Range("S51").Value = 0
$ws.Range("Q52").Value = 0
$ws.Range("S52").Value = 0
$ws.Range("Q53").Value = 0
$ws.Range("R53").Value = 0
$ws.Range("S53").Value = 0

# Header row: add labels for new week columns 14, 15, 16 (bold + centered, like the other header cells)
$ws.Range("Q1:S1").NumberFormat = "@"
$ws.Range("Q1").Value = "14"
$ws.Range("R1").Value = "15"
$ws.Range("S1").Value = "16"
$ws.Range("Q1:S1").Font.Bold = $true
$ws.Range("Q1:S1").HorizontalAlignment = -4108

$ws.Range("A1:S53").Columns.AutoFit() | Out-Null
